# The deck shipped with two embedded theme parts:
#   theme1.xml -> "Office Theme" / clrScheme "Office"      (used by the notes master)
#   theme2.xml -> "Integral"     / clrScheme "Red Violet"  (used by the slide master /
#                                                            the presentation's live design)
# The authored edit swaps the two themes' contents, so the slide master (and therefore
# every slide built from it) now renders with the plain "Office" color palette instead
# of the pink/violet "Integral" palette, while the font scheme and effect/format scheme
# (identical between the two themes already) stay as-is.
#
# The only theme that is reachable/editable through the PowerPoint object model here is
# the one actually driving the presentation's design (the slide master's color scheme),
# so we repoint its twelve theme colors at the "Office" palette values that theme1.xml
# used to hold.

function Convert-HexToComRgb {
    param([string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.ColorScheme

# Target palette = the former theme1.xml ("Office Theme" / clrScheme "Office"), in the
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order exposed via Colors(1..12).
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = Convert-HexToComRgb $officeColors[$i]
}
